$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has GENE numbers in column A (with header blank/missing)
# and the real header/data (TAXON, MODEL_CONDITION, GENE labels, etc.) shifted
# one column to the right (B:F). Remove column A entirely so everything shifts
# left by one column, matching the intended layout (A:E).
$ws.Columns("A").Delete()
